$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as text (avoids Excel
# auto-converting numeric-looking strings like "93.09" into numbers).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Range("D2").Value = "42.052.75"
$ws.Range("E2").Value = "  +5.26%  "

$ws.Range("D3").Value = "2.272.65"
$ws.Range("E3").Value = "  +2.78%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue $ws.Range("D5") "302.52"
$ws.Range("E5").Value = "  +3.62%  "

Set-TextValue $ws.Range("D6") "93.09"
$ws.Range("E6").Value = "  +6.86%  "

Set-TextValue $ws.Range("D7") "0.532"
$ws.Range("E7").Value = "  +4.06%  "

$ws.Range("E8").Value = "  +0.01%  "

Set-TextValue $ws.Range("D9") "0.488"
$ws.Range("E9").Value = "  +4.11%  "

Set-TextValue $ws.Range("D10") "32.79"
$ws.Range("E10").Value = "  +7.48%  "

Set-TextValue $ws.Range("D11") "54.64"
$ws.Range("E11").Value = "  +9.52%  "

Set-TextValue $ws.Range("D12") "0.0802"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("E13").Value = "  +3.08%  "

Set-TextValue $ws.Range("D14") "6.71"
$ws.Range("E14").Value = "  +3.38%  "

$ws.Range("D15").Value = "2.620.67"
$ws.Range("E15").Value = "  +2.67%  "

Set-TextValue $ws.Range("D16") "14.21"
$ws.Range("E16").Value = "  +3.38%  "

$ws.Range("D17").Value = "2.273.63"
$ws.Range("E17").Value = "  -0.53%  "

Set-TextValue $ws.Range("D18") "0.757"
$ws.Range("E18").Value = "  +3.47%  "

$ws.Range("D19").Value = "41.932.19"
$ws.Range("E19").Value = "  +5.25%  "

Set-TextValue $ws.Range("D20") "12.30"
$ws.Range("E20").Value = "  +9.19%  "

$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +2.92%  "

Set-TextValue $ws.Range("D22") "5.95"
$ws.Range("E22").Value = "  +3.38%  "

Set-TextValue $ws.Range("D23") "67.35"
$ws.Range("E23").Value = "  +2.82%  "

Set-TextValue $ws.Range("D24") "242.87"
$ws.Range("E24").Value = "  +2.24%  "

Set-TextValue $ws.Range("D25") "2.58"
$ws.Range("E25").Value = "  +5.75%  "

$ws.Range("E26").Value = "  -0.14%  "

Set-TextValue $ws.Range("D27") "1.92"
$ws.Range("E27").Value = "  +4.85%  "

Set-TextValue $ws.Range("D28") "23.98"
$ws.Range("E28").Value = "  +2.53%  "

Set-TextValue $ws.Range("D29") "2.19"
$ws.Range("E29").Value = "  +6.19%  "

Set-TextValue $ws.Range("D30") "9.73"
$ws.Range("E30").Value = "  +5.36%  "

Set-TextValue $ws.Range("D31") "34.14"
$ws.Range("E31").Value = "  +7.74%  "

Set-TextValue $ws.Range("D32") "157.91"
$ws.Range("E32").Value = "  +0.51%  "

Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("E34").Value = "  +4.80%  "

Set-TextValue $ws.Range("D35") "0.0744"
$ws.Range("E35").Value = "  +4.85%  "

Set-TextValue $ws.Range("D36") "3.09"
$ws.Range("E36").Value = "  +5.93%  "

$ws.Range("E37").Value = "  +3.51%  "

$ws.Range("E38").Value = "  +6.35%  "

$ws.Range("E41").Value = "  +5.75%  "

Set-TextValue $ws.Range("D42") "3.97"
$ws.Range("E42").Value = "  +5.97%  "

Set-TextValue $ws.Range("D43") "20.23"
$ws.Range("E43").Value = "  +14.08%  "

$ws.Range("D44").Value = "2.051.37"
$ws.Range("E44").Value = "  -3.44%  "

Set-TextValue $ws.Range("D45") "0.0280"
$ws.Range("E45").Value = "  +4.31%  "

Set-TextValue $ws.Range("D46") "10.06"
$ws.Range("E46").Value = "  +1.76%  "

Set-TextValue $ws.Range("D47") "2.91"
$ws.Range("E47").Value = "  +8.11%  "

Set-TextValue $ws.Range("D48") "1.99"
$ws.Range("E48").Value = "  -4.43%  "

$ws.Range("D49").Value = "2.492.23"
$ws.Range("E49").Value = "  +2.87%  "

Set-TextValue $ws.Range("D50") "1.52"
$ws.Range("E50").Value = "  +1.84%  "

$ws.Range("E51").Value = "  +4.67%  "

# Swap Celestia (row 39) and Stellar (row 40) data
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D39") "0.116"
$ws.Range("E39").Value = "  +3.78%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D40") "16.58"
$ws.Range("E40").Value = "  +8.10%  "
